$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.403.34"
$ws.Range("E2").Value = "  +3.92%  "
$ws.Range("D3").Value = "2.622.06"
$ws.Range("E3").Value = "  +4.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.174"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.62%  "
$ws.Range("D10").Value = "2.620.21"
$ws.Range("E10").Value = "  +4.18%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "3.104.75"
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("E15").Value = "  +3.82%  "
$ws.Range("D16").Value = "72.213.70"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "2.617.87"
$ws.Range("E18").Value = "  +4.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "381.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +18.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.14%  "
$ws.Range("D28").Value = "2.752.99"
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "0.0₃0953"
$ws.Range("E30").Value = "  +6.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "517.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("E33").Value = "  +6.82%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +6.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.111"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.20%  "
$ws.Range("E41").Value = "  +6.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.49%  "
$ws.Range("E44").Value = "  +8.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("E49").Value = "  +5.04%  "
$ws.Range("E50").Value = "  +7.25%  "
$ws.Range("E51").Value = "  +3.67%  "
